$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E2").Value = 16.465
$ws.Range("C3").Value = -12.903
$ws.Range("C14").Value = -12.139
$ws.Range("C21").Value = -12.63
$ws.Range("C23").Value = -12.223
$ws.Range("C25").Value = -11.988
$ws.Range("E25").Value = 16.98
$ws.Range("C26").Value = -12.937
$ws.Range("E27").Value = 16.19
$ws.Range("C29").Value = -12.041
$ws.Range("E31").Value = 16.6
$ws.Range("E39").Value = 16.47
$ws.Range("E48").Value = 17.314
$ws.Range("E51").Value = 16.878
$ws.Range("E52").Value = 16.848
$ws.Range("C53").Value = -11.515
$ws.Range("E55").Value = 16.491
$ws.Range("E56").Value = 16.4
$ws.Range("C57").Value = -13.794
$ws.Range("E57").Value = 16.572
$ws.Range("C59").Value = -13.155
$ws.Range("C69").Value = -10.697
$ws.Range("E73").Value = 16.778
$ws.Range("C79").Value = -12.083
$ws.Range("C83").Value = -13.169
$ws.Range("E89").Value = 17.362
$ws.Range("E90").Value = 16.537
$ws.Range("C91").Value = -10.595
$ws.Range("E92").Value = 17.779
$ws.Range("C93").Value = -11.766
